{"js": "const REPLACEMENTS = [\n  [\"2023-05-30 Tuesday\", \"2023-05-31 Wednesday\"],\n  [\"34-20=14\", \"1+32=33\"],\n  [\"46-29=17\", \"98-84=14\"],\n  [\"24+64=88\", \"47+46=93\"],\n  [\"29+27=56\", \"27+35=62\"],\n  [\"24+14=38\", \"74-25=49\"],\n  [\"17+63=80\", \"68-28=40\"],\n  [\"83-29=54\", \"6-5=1\"],\n  [\"74+17=91\", \"72-32=40\"],\n  [\"45+24=69\", \"58+14=72\"],\n  [\"75-22=53\", \"85-65=20\"],\n  [\"30+57=87\", \"55-32=23\"],\n  [\"66-53=13\", \"48+43=91\"],\n  [\"71-41=30\", \"72+20=92\"],\n  [\"28+58=86\", \"25+58=83\"],\n  [\"44+15=59\", \"28+31=59\"],\n  [\"63+7=70\", \"17+45=62\"],\n  [\"57+15=72\", \"78-17=61\"],\n  [\"86+1=87\", \"61-3=58\"],\n  [\"77-17=60\", \"64-41=23\"],\n  [\"41+55=96\", \"82-52=30\"],\n  [\"45-11=34\", \"39+15=54\"],\n  [\"23-20=3\", \"5+27=32\"],\n  [\"47+44=91\", \"21+44=65\"],\n  [\"0+64=64\", \"3+9=12\"],\n  [\"64+2=66\", \"94-49=45\"],\n  [\"57+20=77\", \"81-43=38\"],\n  [\"83-55=28\", \"39+58=97\"],\n  [\"98-57=41\", \"77-66=11\"],\n  [\"40-3=37\", \"48-10=38\"],\n  [\"95-89=6\", \"79-61=18\"],\n  [\"78-8=70\", \"74-0=74\"],\n  [\"16-11=5\", \"25+51=76\"],\n  [\"83-0=83\", \"97-0=97\"],\n  [\"68-56=12\", \"85+2=87\"],\n  [\"29+42=71\", \"66-5=61\"],\n  [\"71+21=92\", \"15+25=40\"],\n  [\"73-3=70\", \"5+13=18\"],\n  [\"39-35=4\", \"43+16=59\"],\n  [\"74+9=83\", \"54+7=61\"],\n  [\"34-0=34\", \"2-0=2\"],\n  [\"2+23=25\", \"68-47=21\"],\n  [\"81-59=22\", \"72-2=70\"],\n  [\"22+65=87\", \"91+4=95\"],\n  [\"21+32=53\", \"83-52=31\"],\n  [\"83-58=25\", \"87-55=32\"],\n  [\"74-31=43\", \"8+26=34\"],\n  [\"56+29=85\", \"34+2=36\"],\n  [\"45-4=41\", \"66+21=87\"],\n  [\"52-51=1\", \"6+75=81\"],\n  [\"0+87=87\", \"40+21=61\"],\n  [\"84-1=83\", \"77-21=56\"],\n  [\"47+34=81\", \"10+79=89\"],\n  [\"58-19=39\", \"57-6=51\"],\n  [\"28-8=20\", \"83-64=19\"],\n  [\"17+46=63\", \"25+38=63\"],\n  [\"97-1=96\", \"71+17=88\"],\n  [\"67-18=49\", \"81-58=23\"],\n  [\"4+78=82\", \"35-1=34\"],\n  [\"63+6=69\", \"10+36=46\"],\n  [\"16+13=29\", \"56-54=2\"],\n  [\"82-33=49\", \"71-10=61\"],\n  [\"67+9=76\", \"52+17=69\"],\n  [\"25+47=72\", \"31+67=98\"],\n  [\"69-34=35\", \"72-17=55\"],\n  [\"49+29=78\", \"46-34=12\"],\n  [\"27+47=74\", \"69-60=9\"],\n  [\"71-25=46\", \"1+51=52\"],\n  [\"69-15=54\", \"22-8=14\"],\n  [\"60-31=29\", \"72-55=17\"],\n  [\"17+77=94\", \"42-18=24\"],\n  [\"17+35=52\", \"70+20=90\"],\n  [\"81-76=5\", \"73-16=57\"],\n  [\"19+6=25\", \"94-56=38\"],\n  [\"97-11=86\", \"1+42=43\"],\n  [\"79-47=32\", \"71+20=91\"],\n  [\"26+28=54\", \"16-6=10\"],\n  [\"95-58=37\", \"32+57=89\"],\n  [\"91-2=89\", \"21+78=99\"],\n  [\"48+3=51\", \"26+69=95\"],\n  [\"61-16=45\", \"39+44=83\"],\n  [\"90+1=91\", \"71-43=28\"],\n  [\"14+43=57\", \"91-78=13\"],\n  [\"29-16=13\", \"32-5=27\"],\n  [\"79-5=74\", \"89-40=49\"],\n  [\"25+70=95\", \"44+42=86\"],\n  [\"60-23=37\", \"84-5=79\"],\n  [\"61+27=88\", \"22-1=21\"],\n  [\"0+41=41\", \"50-25=25\"],\n  [\"90-43=47\", \"85+6=91\"],\n  [\"33-2=31\", \"19+63=82\"],\n  [\"93-59=34\", \"77-2=75\"],\n  [\"35-9=26\", \"56+36=92\"],\n  [\"96-25=71\", \"73-43=30\"],\n  [\"44+3=47\", \"10+32=42\"],\n  [\"70+9=79\", \"30+18=48\"],\n  [\"12+27=39\", \"97-82=15\"],\n  [\"5+85=90\", \"73-71=2\"],\n  [\"23+71=94\", \"54+25=79\"],\n  [\"78-14=64\", \"76-21=55\"],\n  [\"3+36=39\", \"50-0=50\"],\n];\n\nfor (const [oldText, newText] of REPLACEMENTS) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for \\\"\" + oldText + \"\\\", found \" + results.items.length);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2023-05-30 Tuesday\", \"2023-05-31 Wednesday\"),\n  @(\"34-20=14\", \"1+32=33\"),\n  @(\"46-29=17\", \"98-84=14\"),\n  @(\"24+64=88\", \"47+46=93\"),\n  @(\"29+27=56\", \"27+35=62\"),\n  @(\"24+14=38\", \"74-25=49\"),\n  @(\"17+63=80\", \"68-28=40\"),\n  @(\"83-29=54\", \"6-5=1\"),\n  @(\"74+17=91\", \"72-32=40\"),\n  @(\"45+24=69\", \"58+14=72\"),\n  @(\"75-22=53\", \"85-65=20\"),\n  @(\"30+57=87\", \"55-32=23\"),\n  @(\"66-53=13\", \"48+43=91\"),\n  @(\"71-41=30\", \"72+20=92\"),\n  @(\"28+58=86\", \"25+58=83\"),\n  @(\"44+15=59\", \"28+31=59\"),\n  @(\"63+7=70\", \"17+45=62\"),\n  @(\"57+15=72\", \"78-17=61\"),\n  @(\"86+1=87\", \"61-3=58\"),\n  @(\"77-17=60\", \"64-41=23\"),\n  @(\"41+55=96\", \"82-52=30\"),\n  @(\"45-11=34\", \"39+15=54\"),\n  @(\"23-20=3\", \"5+27=32\"),\n  @(\"47+44=91\", \"21+44=65\"),\n  @(\"0+64=64\", \"3+9=12\"),\n  @(\"64+2=66\", \"94-49=45\"),\n  @(\"57+20=77\", \"81-43=38\"),\n  @(\"83-55=28\", \"39+58=97\"),\n  @(\"98-57=41\", \"77-66=11\"),\n  @(\"40-3=37\", \"48-10=38\"),\n  @(\"95-89=6\", \"79-61=18\"),\n  @(\"78-8=70\", \"74-0=74\"),\n  @(\"16-11=5\", \"25+51=76\"),\n  @(\"83-0=83\", \"97-0=97\"),\n  @(\"68-56=12\", \"85+2=87\"),\n  @(\"29+42=71\", \"66-5=61\"),\n  @(\"71+21=92\", \"15+25=40\"),\n  @(\"73-3=70\", \"5+13=18\"),\n  @(\"39-35=4\", \"43+16=59\"),\n  @(\"74+9=83\", \"54+7=61\"),\n  @(\"34-0=34\", \"2-0=2\"),\n  @(\"2+23=25\", \"68-47=21\"),\n  @(\"81-59=22\", \"72-2=70\"),\n  @(\"22+65=87\", \"91+4=95\"),\n  @(\"21+32=53\", \"83-52=31\"),\n  @(\"83-58=25\", \"87-55=32\"),\n  @(\"74-31=43\", \"8+26=34\"),\n  @(\"56+29=85\", \"34+2=36\"),\n  @(\"45-4=41\", \"66+21=87\"),\n  @(\"52-51=1\", \"6+75=81\"),\n  @(\"0+87=87\", \"40+21=61\"),\n  @(\"84-1=83\", \"77-21=56\"),\n  @(\"47+34=81\", \"10+79=89\"),\n  @(\"58-19=39\", \"57-6=51\"),\n  @(\"28-8=20\", \"83-64=19\"),\n  @(\"17+46=63\", \"25+38=63\"),\n  @(\"97-1=96\", \"71+17=88\"),\n  @(\"67-18=49\", \"81-58=23\"),\n  @(\"4+78=82\", \"35-1=34\"),\n  @(\"63+6=69\", \"10+36=46\"),\n  @(\"16+13=29\", \"56-54=2\"),\n  @(\"82-33=49\", \"71-10=61\"),\n  @(\"67+9=76\", \"52+17=69\"),\n  @(\"25+47=72\", \"31+67=98\"),\n  @(\"69-34=35\", \"72-17=55\"),\n  @(\"49+29=78\", \"46-34=12\"),\n  @(\"27+47=74\", \"69-60=9\"),\n  @(\"71-25=46\", \"1+51=52\"),\n  @(\"69-15=54\", \"22-8=14\"),\n  @(\"60-31=29\", \"72-55=17\"),\n  @(\"17+77=94\", \"42-18=24\"),\n  @(\"17+35=52\", \"70+20=90\"),\n  @(\"81-76=5\", \"73-16=57\"),\n  @(\"19+6=25\", \"94-56=38\"),\n  @(\"97-11=86\", \"1+42=43\"),\n  @(\"79-47=32\", \"71+20=91\"),\n  @(\"26+28=54\", \"16-6=10\"),\n  @(\"95-58=37\", \"32+57=89\"),\n  @(\"91-2=89\", \"21+78=99\"),\n  @(\"48+3=51\", \"26+69=95\"),\n  @(\"61-16=45\", \"39+44=83\"),\n  @(\"90+1=91\", \"71-43=28\"),\n  @(\"14+43=57\", \"91-78=13\"),\n  @(\"29-16=13\", \"32-5=27\"),\n  @(\"79-5=74\", \"89-40=49\"),\n  @(\"25+70=95\", \"44+42=86\"),\n  @(\"60-23=37\", \"84-5=79\"),\n  @(\"61+27=88\", \"22-1=21\"),\n  @(\"0+41=41\", \"50-25=25\"),\n  @(\"90-43=47\", \"85+6=91\"),\n  @(\"33-2=31\", \"19+63=82\"),\n  @(\"93-59=34\", \"77-2=75\"),\n  @(\"35-9=26\", \"56+36=92\"),\n  @(\"96-25=71\", \"73-43=30\"),\n  @(\"44+3=47\", \"10+32=42\"),\n  @(\"70+9=79\", \"30+18=48\"),\n  @(\"12+27=39\", \"97-82=15\"),\n  @(\"5+85=90\", \"73-71=2\"),\n  @(\"23+71=94\", \"54+25=79\"),\n  @(\"78-14=64\", \"76-21=55\"),\n  @(\"3+36=39\", \"50-0=50\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Could not find text: $oldText\"\n  }\n}\n\n"}
